$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''29.674.51'
$ws.Range("E2").Value = '  +0.63%  '

# Row 3
$ws.Range("D3").Value = '''1.615.99'
$ws.Range("E3").Value = '  +0.72%  '

# Row 4
$ws.Range("E4").Value = '  -0.48%  '

# Row 5
$ws.Range("D5").Value = '''212.69'
$ws.Range("E5").Value = '  +0.15%  '

# Row 6
$ws.Range("E6").Value = '  +0.48%  '

# Row 7
$ws.Range("D7").Value = '''0.992'
$ws.Range("E7").Value = '  -0.52%  '

# Row 8
$ws.Range("D8").Value = '''28.95'
$ws.Range("E8").Value = '  +9.16%  '

# Row 9
$ws.Range("E9").Value = '  +3.06%  '

# Row 10
$ws.Range("D10").Value = '''0.0610'
$ws.Range("E10").Value = '  +1.85%  '

# Row 11
$ws.Range("E11").Value = '  +0.00%  '

# Row 12
$ws.Range("D12").Value = '''1.848.61'
$ws.Range("E12").Value = '  +0.73%  '

# Row 13
$ws.Range("D13").Value = '''1.618.94'
$ws.Range("E13").Value = '  +1.10%  '

# Row 14
$ws.Range("D14").Value = '''0.568'
$ws.Range("E14").Value = '  +6.66%  '

# Row 15
$ws.Range("D15").Value = '''3.87'
$ws.Range("E15").Value = '  +4.87%  '

# Row 16
$ws.Range("D16").Value = '''29.676.80'
$ws.Range("E16").Value = '  +0.54%  '

# Row 17
$ws.Range("D17").Value = '''8.96'
$ws.Range("E17").Value = '  +17.32%  '

# Row 18
$ws.Range("D18").Value = '''64.24'
$ws.Range("E18").Value = '  +1.39%  '

# Row 19
$ws.Range("D19").Value = '''242.38'
$ws.Range("E19").Value = '  +0.14%  '

# Row 20
$ws.Range("D20").Value = '''0.0₃0712'
$ws.Range("E20").Value = '  +3.18%  '

# Row 21
$ws.Range("E21").Value = '  -0.35%  '

# Row 22
$ws.Range("E22").Value = '  +3.19%  '

# Row 23
$ws.Range("D23").Value = '''9.71'
$ws.Range("E23").Value = '  +6.22%  '

# Row 24
$ws.Range("D24").Value = '''2.12'
$ws.Range("E24").Value = '  +1.36%  '

# Row 25
$ws.Range("D25").Value = '''156.71'
$ws.Range("E25").Value = '  +1.48%  '

# Row 26
$ws.Range("E26").Value = '  +2.22%  '

# Row 27
$ws.Range("D27").Value = '''0.111'
$ws.Range("E27").Value = '  +2.09%  '

# Row 28
$ws.Range("D28").Value = '''6.60'
$ws.Range("E28").Value = '  +3.56%  '

# Row 29
$ws.Range("D29").Value = '''0.994'
$ws.Range("E29").Value = '  -0.47%  '

# Row 30
$ws.Range("D30").Value = '''0.0488'
$ws.Range("E30").Value = '  +3.30%  '

# Row 31
$ws.Range("D31").Value = '''3.32'
$ws.Range("E31").Value = '  +3.36%  '

# Row 32
$ws.Range("E32").Value = '  +1.23%  '

# Row 33
$ws.Range("D33").Value = '''3.22'
$ws.Range("E33").Value = '  +3.69%  '

# Row 34
$ws.Range("D34").Value = '''1.435.17'
$ws.Range("E34").Value = '  +1.30%  '

# Row 35
$ws.Range("E35").Value = '  +6.92%  '

# Row 36
$ws.Range("D36").Value = '''1.04'
$ws.Range("E36").Value = '  +1.65%  '

# Row 37
$ws.Range("E37").Value = '  +2.87%  '

# Row 38
$ws.Range("D38").Value = '''2.29'
$ws.Range("E38").Value = '  -0.57%  '

# Row 39
$ws.Range("D39").Value = '''0.0171'
$ws.Range("E39").Value = '  +3.67%  '

# Row 40
$ws.Range("E40").Value = '  +3.84%  '

# Row 41
$ws.Range("D41").Value = '''0.0503'
$ws.Range("E41").Value = '  +3.78%  '

# Row 42
$ws.Range("D42").Value = '''0.829'
$ws.Range("E42").Value = '  +4.53%  '

# Row 43
$ws.Range("E43").Value = '  +0.85%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''69.74'
$ws.Range("E44").Value = '  +6.39%  '

# Row 45
$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").Value = '''53.83'
$ws.Range("E45").Value = '  +0.17%  '

# Row 46
$ws.Range("E46").Value = '  -0.47%  '

# Row 47
$ws.Range("E47").Value = '  +19.77%  '

# Row 48
$ws.Range("E48").Value = '  +3.22%  '

# Row 49
$ws.Range("D49").Value = '''1.757.53'
$ws.Range("E49").Value = '  +0.59%  '

# Row 50
$ws.Range("D50").Value = '''88.15'
$ws.Range("E50").Value = '  +1.79%  '

# Row 51
$ws.Range("E51").Value = '  -0.56%  '
